# Apply the "updated googltrans and requirements" edit:
# Fill in the missing "d" markers in column H of the "algorithms" sheet
# (shared string "d") for every data row that did not already have one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algorithms")

$rows = @(13,15,17,18,19,21,22,23,24,25,26,27,28,29,30,31,32,33,35,36,37,38,39,40,41,42,43,44,46,47,48,49,50,51,52,53,54,56,57,58,59,60,61,62,63,65,66,67,68,69,70,72,75,77,79,81,83,84,85,86,88,89,103,106,108,121,124,126,130)

foreach ($r in $rows) {
    $ws.Range("H$r").Value = "d"
}

# Move the selection to roughly match where the author ended up working.
$ws.Range("J151").Select()
